# Updating strHisCl1 plots to include a comparison between histamine+ iL3s
# and unstimulated wild-type iL3s.
#
# Appends 15 new data rows (194-208) to Sheet1: Experiment="Unstimulated",
# ThermoMode="PT", Phenotype="Unstimulated", with Distance Ratio (E) and
# Speed (F) values, formatted in a smaller Arial font matching the rest of
# the workbook's data-row styling conventions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{Row=194; E=5.8946;  F=0.25887}
    @{Row=195; E=11.373;  F=0.1404}
    @{Row=196; E=6.754;   F=0.14378}
    @{Row=197; E=4.4068;  F=0.14842}
    @{Row=198; E=6.7715;  F=0.14435}
    @{Row=199; E=6.979;   F=0.19669}
    @{Row=200; E=4.3244;  F=0.094409}
    @{Row=201; E=4.9255;  F=0.21005}
    @{Row=202; E=7.2979;  F=0.055652}
    @{Row=203; E=8.4199;  F=0.12308}
    @{Row=204; E=7.5479;  F=0.26831}
    @{Row=205; E=9.7295;  F=0.22291}
    @{Row=206; E=6.2157;  F=0.28038}
    @{Row=207; E=18.493;  F=0.17806}
    @{Row=208; E=12.8;    F=0.067673}
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = "Unstimulated"
    $ws.Cells.Item($r, 2).Value = "PT"
    $ws.Cells.Item($r, 3).Value = "Unstimulated"
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}

# Apply the small Arial fonts used for the Distance Ratio / Speed columns to
# the new rows -- set once on the first row, then fan the resulting format
# out to the rest of the block (avoids re-deriving the font on every cell).
$eDonor = $ws.Cells.Item(194, 5)
$eDonor.Font.Name = "Arial"
$eDonor.Font.Size = 9
$eDonor.Copy()
$ws.Range("E195:E208").PasteSpecial(-4122)

$fDonor = $ws.Cells.Item(194, 6)
$fDonor.Font.Name = "Arial"
$fDonor.Font.Size = 8
$fDonor.Copy()
$ws.Range("F195:F208").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New column C ("Phenotype" continuation) needs a bit more width now that it
# holds "Unstimulated" entries.
$ws.Columns.Item(3).ColumnWidth = 13.3

# Leave the selection where the editor ended up after adding the new block.
[void]$ws.Range("A202").Select()
